# Attendance.xlsx -- "new sheet created for each month"
#
# Before:
#   Tab1 "Sheet1"          -> attendance log (Roll number, Name, 2023-03-12 in/out times)  [active]
#   Tab2 "Student details" -> roster (Roll number, Name, Phone number, Parent name, Parent phone number)
#
# After:
#   Tab1 "Student details" -> roster, contact columns switched from phone numbers to emails
#   Tab2 "March"            -> the old "Sheet1" attendance log, renamed for the month and
#                               now the active/selected tab, with a fresh attendance entry
#
# NOTE: worksheet object references returned by Worksheets.Item(...) are positional, so after
# any .Move() call we must re-fetch sheets by name before touching them again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Re-order the tabs: "Student details" becomes the first tab, "Sheet1" stays
#    second (about to be renamed to "March"), then rename & activate it.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Student details").Move($wb.Worksheets.Item("Sheet1"))
$wb.Worksheets.Item("Sheet1").Name = "March"

$wsStudent = $wb.Worksheets.Item("Student details")
$wsMarch   = $wb.Worksheets.Item("March")

$wsMarch.Activate()

# ---------------------------------------------------------------------------
# 2) "Student details" sheet: switch the contact columns from phone numbers to
#    emails, and make sure the roll number is stored/formatted as a number.
# ---------------------------------------------------------------------------
$wsStudent.Range("C1").Value = "Email"
$wsStudent.Range("E1").Value = "Parent email"

$wsStudent.Range("C2").Value = "piyushchugeja@gmail.com"
$wsStudent.Range("E2").Value = "muskan.chugeja@gmail.com"

$wsStudent.Range("A2").Value = 72

$wsStudent.Range("A1:A2").NumberFormat = "#,##0"
$wsStudent.Range("A1:A2").HorizontalAlignment = -4152
$wsStudent.Range("A1:A2").Font.Color = 0

$wsStudent.Rows.Item(1).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 3) "March" sheet (the former "Sheet1"): rename the roll-number header,
#    record a fresh in/out attendance entry, and tidy up the roll-number
#    cell formatting to match the rest of the workbook.
# ---------------------------------------------------------------------------
$wsMarch.Range("A1").Value = "Roll no"
$wsMarch.Range("C2").Value = "In-time: 23:57:54 " + [char]10 + "Out-time: 23:57:56"

$wsMarch.Range("A2").NumberFormat = "#,##0"
$wsMarch.Range("A2").HorizontalAlignment = -4152
$wsMarch.Range("A2").Font.Color = 0

$wsMarch.Rows.Item(2).RowHeight = 31.5
$wsMarch.Columns.Item(3).ColumnWidth = 17
